$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("PayNowCC_27").Range("B2").Value = "Thu Jun 05 11:57:14 IST 2025"
$wb.Worksheets.Item("PayNowSCFCC_27").Range("B2").Value = "Thu Jun 05 12:00:19 IST 2025"
$wb.Worksheets.Item("PayNowDCFCC_27").Range("B2").Value = "Thu Jun 05 11:53:15 IST 2025"
$wb.Worksheets.Item("CCDeferredCC_27").Range("B2").Value = "Wed Jun 04 12:55:00 IST 2025"
$wb.Worksheets.Item("CMCAutopayCC_27").Range("B2").Value = "Thu Jun 05 11:45:24 IST 2025"
